$d = $word.ActiveDocument

$replacements = @(
    ,@(" Introdução;2) Modos de ", " Introdução;^l2) Modos de ")
    ,@("ia de calor;3) Condução;", "ia de calor;^l3) Condução;")
    ,@("3) Condução;4) Equação d", "3) Condução;^l4) Equação d")
    ,@("da condução;5) Superfíci", "da condução;^l5) Superfíci")
    ,@("s (aletas); 6) Coeficien", "s (aletas); ^l6) Coeficien")
    ,@("o empírico);7) Análise t", "o empírico);^l7) Análise t")
    ,@(" transiente;8) Projeto d", " transiente;^l8) Projeto d")
    ,@(" de energia;2) Modos de ", " de energia;^l2) Modos de ")
    ,@(" e radiação;3) Condução:", " e radiação;^l3) Condução:")
    ,@(" e esférica;4) Equação d", " e esférica;^l4) Equação d")
    ,@(" transiente;5) Superfíci", " transiente;^l5) Superfíci")
    ,@("superfície; 6) Coeficien", "superfície; ^l6) Coeficien")
    ,@("e convecção;7) Análise t", "e convecção;^l7) Análise t")
    ,@("os e ábacos;8) Projeto d", "os e ábacos;^l8) Projeto d")
    ,@(" of Energy; 2) Heat tran", " of Energy; ^l2) Heat tran")
    ,@("d radiation;3) Heat Cond", "d radiation;^l3) Heat Cond")
    ,@("geometries; 4) Different", "geometries; ^l4) Different")
    ,@("te heat conduction; 5) Extended surfaces", "te heat conduction; ^l5) Extended surfaces")
    ,@("iciency fin;6) Convectiv", "iciency fin;^l6) Convectiv")
    ,@("rrelations; 7) Transient", "rrelations; ^l7) Transient")
    ,@("d abacuses; 8) Heat exch", "d abacuses; ^l8) Heat exch")
    ,@(": LTC. 2013.2)KREITH, Fr", ": LTC. 2013.^l2)KREITH, Fr")
    ,@("neira. 2014.3) ÖZISIC, M", "neira. 2014.^l3) ÖZISIC, M")
    ,@("oogan. 1990.4) HOLMAN, J", "oogan. 1990.^l4) HOLMAN, J")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $find
    }
}

Write-Host "Done"
